$wb = $excel.ActiveWorkbook

# --- Sheet 1: "LoopFilter LPF" ---
$ws1 = $wb.Worksheets.Item("LoopFilter LPF")
$ws1.Activate()
$ws1.Range("B3").Value = 7200
$ws1.PageSetup.Orientation = 1
$ws1.Range("B21").Select() | Out-Null

# --- Sheet 2: "Branch LPF" ---
$ws2 = $wb.Worksheets.Item("Branch LPF")
$ws2.Activate()
$ws2.Range("B3").Value = 7200
$ws2.Range("B18").Value = 2
$ws2.PageSetup.Orientation = 1
$ws2.Range("B24").Select() | Out-Null

# Re-activate sheet 1 so it remains the tab shown when opened (tabSelected)
$ws1.Activate()

$wb.Save() | Out-Null
